$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the duplicate fastq row: select the entire row 5 (as a user would by
# clicking the row header) then delete it, which shifts all subsequent rows up by one.
$ws.Rows.Item(5).Select() | Out-Null
$ws.Rows.Item(5).Delete()
